# "updated main GSC export data"
#
# The "Chart" sheet holds a rolling 90-day window of GSC export data in
# A2:C91 (Date, Non-HTTPS URLs, HTTPS URLs). This update rolls the window
# forward by one day:
#   - the oldest date (2025-10-31, row 2) drops off
#   - every remaining row's Date/Non-HTTPS/HTTPS values shift up one row
#   - a new row is appended for the next day (2026-01-29) with its own
#     HTTPS URL page count (27) and Non-HTTPS count (0, same as every
#     other row)
#
# Column A holds the date as literal text (not a real Excel date), so we
# have to defeat Excel's automatic "this text looks like a date" numeric
# coercion: we enter it with a leading apostrophe (forces text) and then
# call ClearFormats() to drop the quote-prefix cell style back to the
# sheet's default, matching the original formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$firstRow = 2
$lastRow = 91

# Snapshot the current Date (col A), Non-HTTPS (col B) and HTTPS (col C)
# columns before overwriting anything, since the shift reads row r+1 while
# writing row r.
$dateVals = @{}
$nonHttpsVals = @{}
$httpsVals = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $dateVals[$r] = $ws.Cells.Item($r, 1).Value2
    $nonHttpsVals[$r] = $ws.Cells.Item($r, 2).Value2
    $httpsVals[$r] = $ws.Cells.Item($r, 3).Value2
}

# Shift rows 3..lastRow up into rows 2..(lastRow-1).
for ($r = $firstRow; $r -le $lastRow - 1; $r++) {
    $ws.Cells.Item($r, 1).Value = "'" + $dateVals[$r + 1]
    $ws.Cells.Item($r, 1).ClearFormats()
    $ws.Cells.Item($r, 2).Value = $nonHttpsVals[$r + 1]
    $ws.Cells.Item($r, 3).Value = $httpsVals[$r + 1]
}

# New last row: the next calendar day after the old last date, with the new
# HTTPS-page count from the commit's data (Non-HTTPS count stays 0,
# consistent with every other row).
$dateParts = $dateVals[$lastRow].Split("-")
$oldLastDate = Get-Date -Year ([int]$dateParts[0]) -Month ([int]$dateParts[1]) -Day ([int]$dateParts[2])
$newLastDate = $oldLastDate.AddDays(1)

$ws.Cells.Item($lastRow, 1).Value = "'" + $newLastDate.ToString("yyyy-MM-dd")
$ws.Cells.Item($lastRow, 1).ClearFormats()
$ws.Cells.Item($lastRow, 2).Value = 0
$ws.Cells.Item($lastRow, 3).Value = 27
